$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New "gm" row (row 8), inserted right after the connect-redis row ---
$ws.Cells.Item(8, 1).Value = "gm"
$ws.Cells.Item(8, 3).Value = "npm install gm --save"
$ws.Cells.Item(8, 5).Value = "要预先安装GraphicMagic`nftp://ftp.graphicsmagick.org/pub/GraphicsMagick/windows/`n安装完毕，将gm目录加入环境变量"

# New column E needs extra width to hold the GraphicsMagick note, and the note cell
# wraps its text across the taller row.
$ws.Columns.Item(5).ColumnWidth = 23.43
$ws.Cells.Item(8, 5).WrapText = $true
$ws.Rows.Item(8).RowHeight = 81

# D7 (connect-redis hyperlink) now shares the same named hyperlink style used by the
# other link cells in the sheet (D17/D18), instead of its own ad-hoc duplicate.
$ws.Range("D7").Style = "超链接"

# Match the author's final selection state.
$ws.Range("G13:G15").Select()
